# Update SwaadSutra_Daily_2026-01-19.xlsx
# New order #18 (Radhika Joshi, Pohe x3) arrived -> inserted as the newest row
# at the top of the "Daily Orders" table, pushing the existing orders down by
# one row. Summary + Items Breakdown roll-ups are updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Daily Orders" sheet - insert the new order as row 2.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Daily Orders")

# Stage the new row's values in a scratch row far below the real data so we
# can fix up the "text that looks like a number/date" columns (Phone,
# Collection Date) with a Text number format *before* they are ever
# assigned - this keeps them as real text instead of Excel auto-converting
# them into a number / date serial, and it does so without leaving that
# Text format behind on the final cell (the values get carried over via a
# copy into the freshly-inserted row, which preserves the General/default
# style of that row).
$scratchRow = 20

$ws.Range("A" + $scratchRow).Value = 18
$ws.Range("B" + $scratchRow).Value = "2026-01-19 08:37"
$ws.Range("C" + $scratchRow).Value = "Radhika Joshi"
$ws.Range("D" + $scratchRow).Value = "C 1501"
$ws.Range("E" + $scratchRow).NumberFormat = "@"
$ws.Range("E" + $scratchRow).Value = "9967195227"
$ws.Range("F" + $scratchRow).Value = "Pohe x3"
$ws.Range("G" + $scratchRow).Value = 90
$ws.Range("H" + $scratchRow).Value = "NEW"
$ws.Range("I" + $scratchRow).Value = "PENDING"
$ws.Range("J" + $scratchRow).NumberFormat = "@"
$ws.Range("J" + $scratchRow).Value = "2026-01-20"
$ws.Range("K" + $scratchRow).Value = "08:00"

# Copy the staged row and insert it above row 2 - this shifts the existing
# orders (17, 16, 15) down to rows 3, 4, 5 and drops the copied values into
# the now-empty row 2.
$ws.Range("A" + $scratchRow + ":K" + $scratchRow).Copy()
$ws.Range("A2").EntireRow.Insert()

# Remove the scratch row (it shifted down to row 21 when row 2 was inserted
# above it).
$ws.Range("A21").EntireRow.Delete()

# Notes / Cancel Reason / Feedback are blank for the new order, same as the
# other rows.
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'"
$ws.Range("N2").Value = "'"

# ---------------------------------------------------------------------------
# 2) "Summary" sheet - one more order, +90 revenue.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = 4
$summary.Range("B2").Value = 4
$summary.Range("G2").Value = 375

# ---------------------------------------------------------------------------
# 3) "Items Breakdown" sheet - add the "Pohe" line item above "1 Plate Bhaji".
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Range("A3").EntireRow.Insert()
$items.Range("A3").Value = "Pohe"
$items.Range("B3").Value = 3
$items.Range("C3").Value = 90
